$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing account 005081833 (PEDRO, 100000).
# That record is on worksheet row 5 (row 1 = header "Conta/Nome/Saldo";
# row 2 = CAROLINA; row 3 = BRASFORT; row 4 = GABRIEL; row 5 = PEDRO).
$ws.Rows.Item(5).Delete()
